# Insert a new weekly price entry as row 302 on the (single) active sheet.
# This pushes the existing rows 302-346 down to 303-347, and appends the
# data that used to live at the bottom (row 346) as the new last row (347) -
# i.e. a plain row insert, Excel reflows everything below it automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 302.
$ws.Rows.Item(302).Insert()

# Fill the new row 302 with the new weekly record.
$ws.Cells.Item(302, 1).Value = 3
$ws.Cells.Item(302, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(302, 3).Value = "Coquimbo"
$ws.Cells.Item(302, 4).Value = 44776
$ws.Cells.Item(302, 5).Value = 5
$ws.Cells.Item(302, 6).Value = 100112039
$ws.Cells.Item(302, 7).Value = "Ciboulette"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 120
$ws.Cells.Item(302, 11).Value = 1500
$ws.Cells.Item(302, 12).Value = 1500
$ws.Cells.Item(302, 13).Value = 1500
$ws.Cells.Item(302, 14).Value = "`$/docena de atados"
$ws.Cells.Item(302, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(302, 16).Value = 500
$ws.Cells.Item(302, 17).Value = 3
$ws.Cells.Item(302, 18).Value = "Hortaliza"
